$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Bring this sheet to the front / make it the active tab (matches
# workbookView activeTab="3" and sheetView tabSelected="1" moving here).
[void]$ws.Activate()

# Insert a new blank column before column N (Late/Over Due/heading/Outstanding
# shift one column to the right: N->O, O->P, P->Q). Inherit the column width
# from the column to its left (M), which is what Excel does natively when a
# column is inserted.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $mWidth

# Match the new selection left behind on this sheet.
$ws.Range("S6").Select() | Out-Null
